$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.926.27"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "1.912.62"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.93"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5004"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2998"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("D10").Value = "1.908.72"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.06"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07319"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.42"
$ws.Range("E13").Value = "  +6.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.118"
$ws.Range("E14").Value = "  +5.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6840"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "30.898.79"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008059"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.31"
$ws.Range("E18").Value = "  +4.27%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "2.154.15"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.891"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "183.71"
$ws.Range("E23").Value = "  +35.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.138"
$ws.Range("E24").Value = "  +9.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.391"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.02"
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +12.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.395"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.358"
$ws.Range("E30").Value = "  +4.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09002"
$ws.Range("E31").Value = "  +3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.067"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05272"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7495"
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.670"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01961"
$ws.Range("E37").Value = "  +19.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.740"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.186"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9371"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4414"
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.54"
$ws.Range("E42").Value = "  +5.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.880"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.798"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1350"
$ws.Range("E46").Value = "  +7.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05846"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3944"
$ws.Range("E48").Value = "  +6.50%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.634"
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.33"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.397"
$ws.Range("E51").Value = "  +4.58%  "
